$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsmatrix")

function Copy-Fmt($srcRange, $dstRange) {
    $ws.Range($srcRange).Copy()
    $ws.Range($dstRange).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# 1) Apply per-row cell formatting (style indices) BEFORE writing any values,
#    by copying formats from structurally-identical existing rows. We copy in
#    contiguous column blocks that exactly match which columns should exist
#    for the destination row, so no stray cells get created.
# ---------------------------------------------------------------------------

# Row 55: A:K populated (incl. H formula cell) -> same shape as row 52
Copy-Fmt "A52:K52" "A55:K55"

# Row 56: A:G + I:K populated (no H) -> same shape as row 21
Copy-Fmt "A21:G21" "A56:G56"
Copy-Fmt "I21:K21" "I56:K56"

# Row 57: A:G + I:K + L:M populated (no H) -> A:G/I:K like row 21, L:M like row 49
Copy-Fmt "A21:G21" "A57:G57"
Copy-Fmt "I21:K21" "I57:K57"
Copy-Fmt "L49:M49" "L57:M57"

# Row 58: only G, I, J, K populated
Copy-Fmt "G52" "G58"
Copy-Fmt "I52" "I58"
Copy-Fmt "J52" "J58"
Copy-Fmt "K52" "K58"

# Row 59: same shape as row 21 (A:G + I:K, no H)
Copy-Fmt "A21:G21" "A59:G59"
Copy-Fmt "I21:K21" "I59:K59"

# Row 60: same shape as row 21 (A:G + I:K, no H)
Copy-Fmt "A21:G21" "A60:G60"
Copy-Fmt "I21:K21" "I60:K60"

# Row 61: A:K populated (incl. H formula cell) -> same shape as row 52
Copy-Fmt "A52:K52" "A61:K61"

# Row 62: same shape as row 21 (A:G + I:K, no H)
Copy-Fmt "A21:G21" "A62:G62"
Copy-Fmt "I21:K21" "I62:K62"

# Rows 63 & 64: only I populated
Copy-Fmt "I52" "I63"
Copy-Fmt "I52" "I64"

# Row 65 (new totals row): B, C, D, E, F, G populated -> same shape as old row 60
Copy-Fmt "B60:G60" "B65:G65"

# ---------------------------------------------------------------------------
# 2) Write cell values / formulas for the new rows
# ---------------------------------------------------------------------------

# Row 55
$ws.Range("A55").Value = 15
$ws.Range("B55").Value = "Konzeptuelles Design"
$ws.Range("C55").Value = "[SEMINAR]"
$ws.Range("D55").Value = "Peer Reviewed Expose"
$ws.Range("E55").Value = "Expose der anderen kommentieren"
$ws.Range("F55").Value = 44331
$ws.Range("G55").Value = 44338
$ws.Range("H55").Formula = "=ROUNDUP(((SUM(K55-J55)*24*60/60)/0.25),0)*0.25"
$ws.Range("J55").Value = 0.375
$ws.Range("K55").Value = 0.5

# Row 56
$ws.Range("A56").Value = 14
$ws.Range("B56").Value = "Konzeptuelles Design"
$ws.Range("C56").Value = "[FEATURE]"
$ws.Range("D56").Value = "Umfrage"
$ws.Range("E56").Value = "Umfrage erstellt und versendet"
$ws.Range("F56").Value = 44331
$ws.Range("G56").Value = 44338
$ws.Range("I56").Formula = "=ROUNDUP(((SUM(K56-J56)*24*60/60)/0.25),0)*0.25"
$ws.Range("J56").Value = 0.5
$ws.Range("K56").Value = 0.57291666666666663

# Row 57
$ws.Range("A57").Value = 14
$ws.Range("B57").Value = "Konzeptuelles Design"
$ws.Range("C57").Value = "[TASK]"
$ws.Range("D57").Value = "Umfrage"
$ws.Range("E57").Value = "Umfrage auswerten"
$ws.Range("F57").Value = 44332
$ws.Range("G57").Value = 44338
$ws.Range("I57").Formula = "=ROUNDUP(((SUM(K57-J57)*24*60/60)/0.25),0)*0.25"
$ws.Range("J57").Value = 0.70833333333333337
$ws.Range("K57").Value = 0.79166666666666663
$ws.Range("L57").Formula = "=SUM(H49:I57)"
$ws.Range("M57").Formula = "=SUM(L57+16)"

# Row 58 stays blank (only carries formatting, no values)

# Row 59
$ws.Range("A59").Value = 18
$ws.Range("B59").Value = "Konzeptuelles Design"
$ws.Range("C59").Value = "[FEATURE]"
$ws.Range("D59").Value = "Content Map"
$ws.Range("E59").Value = "Content Map beginnen"
$ws.Range("F59").Value = 44333
$ws.Range("G59").Value = 44338
$ws.Range("I59").Formula = "=ROUNDUP(((SUM(K59-J59)*24*60/60)/0.25),0)*0.25"
$ws.Range("J59").Value = 0.41666666666666669
$ws.Range("K59").Value = 0.54166666666666663

# Row 60
$ws.Range("A60").Value = 18
$ws.Range("B60").Value = "Konzeptuelles Design"
$ws.Range("C60").Value = "[TASK]"
$ws.Range("D60").Value = "Content Map"
$ws.Range("E60").Value = "Content Map iterieren"
$ws.Range("F60").Value = 44333
$ws.Range("G60").Value = 44338
$ws.Range("I60").Formula = "=ROUNDUP(((SUM(K60-J60)*24*60/60)/0.25),0)*0.25"
$ws.Range("J60").Value = 0.625
$ws.Range("K60").Value = 0.70833333333333337

# Row 61
$ws.Range("A61").Value = 16
$ws.Range("B61").Value = "Konzeptuelles Design"
$ws.Range("C61").Value = "[SEMINAR]"
$ws.Range("D61").Value = "Peer reviewed Exposé"
$ws.Range("E61").Value = "Exposé vorstellen"
$ws.Range("F61").Value = 44334
$ws.Range("G61").Value = 44338
$ws.Range("H61").Formula = "=ROUNDUP(((SUM(K61-J61)*24*60/60)/0.25),0)*0.25"
$ws.Range("J61").Value = 0.41666666666666669
$ws.Range("K61").Value = 0.54166666666666663

# Row 62
$ws.Range("A62").Value = 18
$ws.Range("B62").Value = "Konzeptuelles Design"
$ws.Range("C62").Value = "[TASK]"
$ws.Range("D62").Value = "Content Map"
$ws.Range("E62").Value = "Content Map iterieren"
$ws.Range("F62").Value = 44334
$ws.Range("G62").Value = 44338
$ws.Range("I62").Formula = "=ROUNDUP(((SUM(K62-J62)*24*60/60)/0.25),0)*0.25"
$ws.Range("J62").Value = 0.625
$ws.Range("K62").Value = 0.83333333333333337

# Rows 63 & 64 stay blank (only carry formatting, no values)

# Row 65 - totals row (moved down from row 60)
$ws.Range("B65").Value = "Stunden insgesamt"
$ws.Range("C65").Formula = "=SUM(I:I)+SUM(H:H)"
$ws.Range("D65").Value = "Stunden Seminar"
$ws.Range("E65").Formula = "=SUM(H:H)"
$ws.Range("F65").Value = "Stunden Projekt"
$ws.Range("G65").Formula = "=SUM(I:I)"

# ---------------------------------------------------------------------------
# 3) Clear the old totals-row formulas/values that used to live on row 60
#    (B60:G60) since row 60 is now a regular data row - already overwritten
#    above by the Value assignments, nothing further required.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 4) Sheet view: scroll position + active selection
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 38
$ws.Range("K63").Select()

Write-Output "edit applied"
